$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2984.1538
$ws.Range("I98").Value = 2279.5
$ws.Range("J98").Value = 5333
$ws.Range("K98").Value = 2279.5
$ws.Range("L98").Value = 5333
$ws.Range("M98").Value = -781.5
$ws.Range("N98").Value = -8329
$ws.Range("H107").Value = 877.9394
$ws.Range("I107").Value = 1096.0435
$ws.Range("K107").Value = 1096.0435
$ws.Range("M107").Value = 823.9565
$ws.Range("H122").Value = 2984.1538
$ws.Range("I122").Value = 2279.5
$ws.Range("J122").Value = 5333
$ws.Range("K122").Value = 6838.5
$ws.Range("L122").Value = 15999
$ws.Range("M122").Value = -4388.5
$ws.Range("N122").Value = -20899
$ws.Range("H131").Value = 4821.467
$ws.Range("I131").Value = 4028.5
$ws.Range("J131").Value = 7993.3335
$ws.Range("K131").Value = 12085.5
$ws.Range("L131").Value = 23980.0005
$ws.Range("M131").Value = -7045.5
$ws.Range("N131").Value = -34060.00049999999
$ws.Range("H132").Value = 12154.625
$ws.Range("I132").Value = 13676.714
$ws.Range("J132").Value = 1500
$ws.Range("K132").Value = 41030.142
$ws.Range("L132").Value = 4500
$ws.Range("M132").Value = -38500.142
$ws.Range("N132").Value = -9560
$ws.Range("H135").Value = 1349.7646
$ws.Range("I135").Value = 995.6667
$ws.Range("K135").Value = 8961.0003
$ws.Range("M135").Value = -6426.0003
$ws.Range("H137").Value = 114811.56
$ws.Range("I137").Value = 225309.75
$ws.Range("K137").Value = 675929.25
$ws.Range("M137").Value = -673379.25
$ws.Range("H138").Value = 3014.5908
$ws.Range("J138").Value = 3398.3958
$ws.Range("L138").Value = 10195.1874
$ws.Range("N138").Value = -20475.1874

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4526.2754
$ws.Range("I32").Value = 2434.4312
$ws.Range("K32").Value = 2434.4312
$ws.Range("M32").Value = -2147.4312
$ws.Range("H45").Value = 66151.44
$ws.Range("I45").Value = 126691.25
$ws.Range("J45").Value = 5611.625
$ws.Range("K45").Value = 126691.25
$ws.Range("L45").Value = 5611.625
$ws.Range("M45").Value = -126314.25
$ws.Range("N45").Value = -6365.625
$ws.Range("H61").Value = 4296.3335
$ws.Range("I61").Value = 3947
$ws.Range("J61").Value = 4995
$ws.Range("K61").Value = 3947
$ws.Range("L61").Value = 4995
$ws.Range("M61").Value = -3735
$ws.Range("N61").Value = -5419
$ws.Range("H97").Value = 8894.666999999999
$ws.Range("I97").Value = 8894.666999999999
$ws.Range("K97").Value = 8894.666999999999
$ws.Range("M97").Value = -8398.666999999999
$ws.Range("H132").Value = 2050.8635
$ws.Range("I132").Value = 1562.4445
$ws.Range("K132").Value = 4687.333500000001
$ws.Range("M132").Value = -2157.333500000001
$ws.Range("H136").Value = 4296.3335
$ws.Range("I136").Value = 3947
$ws.Range("J136").Value = 4995
$ws.Range("K136").Value = 11841
$ws.Range("L136").Value = 14985
$ws.Range("M136").Value = -9291
$ws.Range("N136").Value = -20085

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3676.8164
$ws.Range("I86").Value = 5459.5
$ws.Range("K86").Value = 5459.5
$ws.Range("M86").Value = -4336.5
$ws.Range("H89").Value = 3676.8164
$ws.Range("I89").Value = 5459.5
$ws.Range("K89").Value = 27297.5
$ws.Range("M89").Value = -21681.5
$ws.Range("H94").Value = 12049.75
$ws.Range("J94").Value = 17900
$ws.Range("L94").Value = 17900
$ws.Range("N94").Value = -18802
$ws.Range("H107").Value = 1000
$ws.Range("I107").Value = 1000
$ws.Range("K107").Value = 1000
$ws.Range("M107").Value = 920
$ws.Range("H126").Value = 15000
$ws.Range("I126").Value = 15000
$ws.Range("K126").Value = 15000
$ws.Range("M126").Value = -10060
$ws.Range("H134").Value = 5090.6665
$ws.Range("I134").Value = 2570.3572
$ws.Range("K134").Value = 7711.071599999999
$ws.Range("M134").Value = -5176.071599999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 47786.25
$ws.Range("I31").Value = 1811.7273
$ws.Range("K31").Value = 1811.7273
$ws.Range("M31").Value = -1516.7273
$ws.Range("H34").Value = 47786.25
$ws.Range("I34").Value = 1811.7273
$ws.Range("K34").Value = 1811.7273
$ws.Range("M34").Value = -1609.7273
$ws.Range("H58").Value = 1961.0385
$ws.Range("I58").Value = 1942.4166
$ws.Range("K58").Value = 1942.4166
$ws.Range("M58").Value = -1739.4166
$ws.Range("H62").Value = 3000
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 3000
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H107").Value = 2379.8823
$ws.Range("I107").Value = 2323.7693
$ws.Range("K107").Value = 2323.7693
$ws.Range("M107").Value = -403.7692999999999
$ws.Range("H132").Value = 62097.066
$ws.Range("I132").Value = 2740.889
$ws.Range("J132").Value = 151131.33
$ws.Range("K132").Value = 8222.667000000001
$ws.Range("L132").Value = 453393.99
$ws.Range("M132").Value = -5692.667000000001
$ws.Range("N132").Value = -458453.99
$ws.Range("H134").Value = 2762.862
$ws.Range("I134").Value = 2109.158
$ws.Range("J134").Value = 4004.9
$ws.Range("K134").Value = 6327.474
$ws.Range("L134").Value = 12014.7
$ws.Range("M134").Value = -3792.474
$ws.Range("N134").Value = -17084.7
$ws.Range("H136").Value = 1961.0385
$ws.Range("I136").Value = 1942.4166
$ws.Range("K136").Value = 5827.2498
$ws.Range("M136").Value = -3277.2498

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1072.7428
$ws.Range("I5").Value = 709.6667
$ws.Range("J5").Value = 1617.3572
$ws.Range("K5").Value = 2129.0001
$ws.Range("L5").Value = 4852.071599999999
$ws.Range("M5").Value = -2017.0001
$ws.Range("N5").Value = -5076.071599999999
$ws.Range("H68").Value = 642
$ws.Range("I68").Value = 599.2857
$ws.Range("J68").Value = 791.5
$ws.Range("K68").Value = 1797.8571
$ws.Range("L68").Value = 2374.5
$ws.Range("M68").Value = -986.8571000000002
$ws.Range("N68").Value = -3996.5
$ws.Range("H71").Value = 642
$ws.Range("I71").Value = 599.2857
$ws.Range("J71").Value = 791.5
$ws.Range("K71").Value = 5393.571300000001
$ws.Range("L71").Value = 7123.5
$ws.Range("M71").Value = -1337.571300000001
$ws.Range("N71").Value = -15235.5
$ws.Range("H92").Value = 368
$ws.Range("I92").Value = 415
$ws.Range("J92").Value = 358.6
$ws.Range("K92").Value = 1245
$ws.Range("L92").Value = 1075.8
$ws.Range("M92").Value = 3
$ws.Range("N92").Value = -3571.8
$ws.Range("H107").Value = 686.2727
$ws.Range("I107").Value = 217.66667
$ws.Range("J107").Value = 1248.6
$ws.Range("K107").Value = 653.00001
$ws.Range("L107").Value = 3745.8
$ws.Range("M107").Value = 1266.99999
$ws.Range("N107").Value = -7585.799999999999
$ws.Range("H112").Value = 120
$ws.Range("J112").Value = 115
$ws.Range("L112").Value = 345
$ws.Range("N112").Value = -2561
$ws.Range("H135").Value = 1072.7428
$ws.Range("I135").Value = 709.6667
$ws.Range("J135").Value = 1617.3572
$ws.Range("K135").Value = 6387.0003
$ws.Range("L135").Value = 14556.2148
$ws.Range("M135").Value = -3852.0003
$ws.Range("N135").Value = -19626.2148
$ws.Range("H139").Value = 3699.2
$ws.Range("I139").Value = 3332.5
$ws.Range("K139").Value = 9997.5
$ws.Range("M139").Value = -4857.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 87561.25
$ws.Range("I102").Value = 3291.6667
$ws.Range("K102").Value = 3291.6667
$ws.Range("M102").Value = -1669.6667
$ws.Range("H126").Value = 3554.3
$ws.Range("I126").Value = 3394.2222
$ws.Range("J126").Value = 4995
$ws.Range("K126").Value = 10182.6666
$ws.Range("L126").Value = 14985
$ws.Range("M126").Value = -7712.6666
$ws.Range("N126").Value = -19925
$ws.Range("H132").Value = 4201.15
$ws.Range("I132").Value = 3338
$ws.Range("K132").Value = 10014
$ws.Range("M132").Value = -7484

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5482.92
$ws.Range("I40").Value = 4760.0586
$ws.Range("K40").Value = 4760.0586
$ws.Range("M40").Value = -4624.0586
$ws.Range("H55").Value = 1721.6333
$ws.Range("I55").Value = 1447.875
$ws.Range("J55").Value = 2034.5
$ws.Range("K55").Value = 1447.875
$ws.Range("L55").Value = 2034.5
$ws.Range("M55").Value = -1274.875
$ws.Range("N55").Value = -2380.5
$ws.Range("H61").Value = 2774.88
$ws.Range("I61").Value = 2838.4167
$ws.Range("J61").Value = 1250
$ws.Range("K61").Value = 2838.4167
$ws.Range("L61").Value = 1250
$ws.Range("M61").Value = -2636.4167
$ws.Range("N61").Value = -1654
$ws.Range("H113").Value = 2774.88
$ws.Range("I113").Value = 2838.4167
$ws.Range("J113").Value = 1250
$ws.Range("K113").Value = 2838.4167
$ws.Range("L113").Value = 1250
$ws.Range("M113").Value = -668.4167000000002
$ws.Range("N113").Value = -5590

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1547.125
$ws.Range("I81").Value = 1479.8334
$ws.Range("K81").Value = 2959.6668
$ws.Range("M81").Value = -1898.6668
$ws.Range("H84").Value = 1547.125
$ws.Range("I84").Value = 1479.8334
$ws.Range("K84").Value = 14798.334
$ws.Range("M84").Value = -9494.333999999999
$ws.Range("H132").Value = 87126.836
$ws.Range("I132").Value = 18534.285
$ws.Range("K132").Value = 55602.855
$ws.Range("M132").Value = -53072.855
$ws.Range("H136").Value = 4702.923
$ws.Range("I136").Value = 3959
$ws.Range("J136").Value = 5167.875
$ws.Range("K136").Value = 11877
$ws.Range("L136").Value = 15503.625
$ws.Range("M136").Value = -9327
$ws.Range("N136").Value = -20603.625
